$d = $word.ActiveDocument
$bm = $d.Bookmarks("__DdeLink__332_701596780")
$bm.Delete
Write-Output ("exists after=" + $d.Bookmarks.Exists("__DdeLink__332_701596780"))
